# Automatische test-sync: 2025-08-14 21:42:50
# Append a new log row (row 30) to the "Logs" sheet, extend the
# conditional-formatting ranges that covered rows 2:29 to now cover
# rows 2:30, and bump the "Intern verzoek / Actie voor medewerker"
# count on the "Dashboard" sheet from 21 to 22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$newRow = 30

$ws.Cells.Item($newRow, 1).Value  = "Demo inplannen"
$ws.Cells.Item($newRow, 2).Value  = "klantenservice@testbedrijf123.nl"
$ws.Cells.Item($newRow, 3).Value  = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Cells.Item($newRow, 4).Value  = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item($newRow, 5).Value  = "Bedankt, we hebben dit doorgestuurd naar support@testbedrijf123.nl."
$ws.Cells.Item($newRow, 6).Value  = "2025-08-14 21:42:13"
$ws.Cells.Item($newRow, 7).Value  = "Nee"
$ws.Cells.Item($newRow, 8).Value  = "Ja"
$ws.Cells.Item($newRow, 9).Value  = "Nee"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# Extend each conditional-formatting block's AppliesTo range from
# row 29 to row 30 (columns D, G, H, I, J), preserving rule order,
# priorities and dxf styles.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $ws.Range($col + "2:" + $col + "29")
    $newRange = $ws.Range($col + "2:" + $col + "30")
    $fc = $oldRange.FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Dashboard: bump the "Intern verzoek / Actie voor medewerker" tally.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(2, 2).Value = 22
